$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update best_score (B) and best_time (D) for quiz rows 1956, 1957, 1958
$ws.Range("B54").Value = 2
$ws.Range("D54").Value = 13

$ws.Range("B55").Value = 2
$ws.Range("D55").Value = 16

$ws.Range("B56").Value = 2
$ws.Range("D56").Value = 35
